# LOB1206.docx edit: reorder/rotate text blocks among fixed paragraph
# "slots" without altering paragraph styles/structure (matches the
# target XML diff exactly: paragraph count, styles and runs stay the
# same; only the w:t contents move between paragraphs).

$d = $word.ActiveDocument

# --- Whole-paragraph text swaps (single-run paragraphs) ---
# Paragraph 6: Objetivos (PT) body -> becomes short program summary (PT)
$d.Paragraphs.Item(6).Range.Text = "Introdução. Formação do solo. Atributos físicos do solo. Classificação do solo. Água do solo. Aula Prática: Descrição de perfil no campo. Aula Prática: Caracterização e métodos de determinação de atributos físicos e hídricos do solo."

# Paragraph 7: Objetivos (EN, italic) -> becomes short program summary (EN, italic)
$d.Paragraphs.Item(7).Range.Text = "Introduction. Soil formation. Soil physical properties. Soil classification. Soil water. Practical class: Profile description in the field. Practical class: Characterization and determination methods of physical and hydraulic properties of the soil."

# Paragraph 9: Docente bullet -> becomes Objetivos (PT) text
$d.Paragraphs.Item(9).Range.Text = "A disciplina tem o objetivo de apresentar ao estudante informações a respeito das características dos solos, particularmente os existentes na região tropical, e dos fenômenos físicos que nele ocorrem, a fim de capacitá-lo a compreender a importância dos fatores pedológicos, físicos e hídricos na preservação do ambiente."

# Paragraph 11: short program summary (PT) -> becomes long program (PT)
$d.Paragraphs.Item(11).Range.Text = "INTRODUÇÃO. Conceitos Básicos. O perfil de solo. Definição e notação de horizontes e camadas. FORMAÇÃO DO SOLO. Fatores e processos de formação. Intemperismo. ATRIBUTOS FÍSICOS DO SOLO. Composição volumétrica, granulometria e textura, estrutura e agregação, cor, porosidade, densidade e compactação, consistência. CLASSIFICAÇÃO DO SOLO. Sistema brasileiro de classificação de solos. Principais atributos morfológicos. Principais Classes de Solos. ÁGUA DO SOLO. Conceito e importância. Constantes de umidade. Potencial total da água do solo e seus componentes. Curva característica da água do solo. Movimento da Água e de solutos no Solo. Aula prática de campo: Descrição de perfil no campo. Aula prática de laboratório: Caracterização e métodos de determinação de atributos físicos e hídricos do solo. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."

# Paragraph 12: short program summary (EN, italic) -> becomes Objetivos (EN, italic) text
$d.Paragraphs.Item(12).Range.Text = "The course aims to introduce the student information about the soil characteristics, particularly in the soil existing in the tropics, and the physical phenomena occurring in it in order to enable the students to understand the importance of pedological, physical and hydric factors in the environment preservation."

# Paragraph 14: long program (PT) -> becomes evaluation "Método" text
$d.Paragraphs.Item(14).Range.Text = "O aluno poderá optar por um dos dois critérios de avaliação para a NF (nota final).  Critério 1: NF = média obtida em todas atividades desenvolvidas, trabalhos e relatórios ao longo do semestre. Critério 2 (alternativo): NF = (P1+P2)/2, sendo P1 e P2 avaliações escritas individuais."

# Paragraph 19: bibliography text -> becomes Docente bullet text
$d.Paragraphs.Item(19).Range.Text = "5840942 - Marco Aurélio Kondracki de Alcântara"

# --- Paragraph 17 (Avaliação bullet, multiple runs) ---
# Runs in order: "Método: " | <method text> | "Critério: " | <criteria text> |
#                "Norma de recuperação: " | <norma text>
# Only the "Método:" run's text and the "Norma de recuperação:" run's text
# change; the "Critério:" run's text is left as-is (it duplicates the
# original text that also used to be the "Norma de recuperação:" text).
$p17 = $d.Paragraphs.Item(17)
$p17Start = $p17.Range.Start
$p17End = $p17.Range.End
$p17Text = $p17.Range.Text

# Replace the "Método:" run's value. Scope the Find to the paragraph range
# but stop before the "Critério:" label so the identical text that remains
# under "Critério:" is not touched.
$criterioIdx = $p17Text.IndexOf("Critério: ")
$metodoScope = $d.Range($p17Start, $p17Start + $criterioIdx)
$ok1 = $metodoScope.Find.Execute(
    "O aluno poderá optar por um dos dois critérios de avaliação para a NF (nota final).  Critério 1: NF = média obtida em todas atividades desenvolvidas, trabalhos e relatórios ao longo do semestre. Critério 2 (alternativo): NF = (P1+P2)/2, sendo P1 e P2 avaliações escritas individuais.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.",
    2)
if (-not $ok1) { throw "Find/replace for the 'Método:' run failed" }

# Replace the "Norma de recuperação:" run's value. Scope the Find to begin
# at the "Norma de recuperação:" label so the earlier identical text under
# "Critério:" is not touched.
$p17Text2 = $p17.Range.Text
$normaIdx = $p17Text2.IndexOf("Norma de recuperação")
$normaScope = $d.Range($p17Start + $normaIdx, $p17.Range.End)
$ok2 = $normaScope.Find.Execute(
    "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bibliografia básica:1. MEURER, E.J. Fundamentos de Química do Solo, 3.ed. Porto Alegre: Editora Evangraf, 2010. 264p.2. ERNANI, P.R. Química do Solo e Disponibilidade de Nutrientes. Lages: Udesc, 1. ed. , 2008. v.1. 230 p.Bibliografia complementar:1. CAMARGO, O.A. de; MONIZ, A.C.; JORGE, J.A.; VALADARES, J.M.A.S. Métodos de analise química, mineralógica e física de solos do Instituto Agronômico de Campinas. Campinas, Instituto Agronômico, 2009. 77 p. (Boletim técnico, 106, Edição revista e atualizada).2. EMBRAPA. Centro Nacional de Pesquisa de Solos (Rio de Janeiro, RJ). Sistema Brasileiro de Classificação de Solos. Brasília: Embrapa Produção da Informação; Rio de Janeiro: Embrapa Solos, 2006. 306p.3. MELO, V.F.; ALLEONI, L.R.F. (Eds.). Química e mineralogia do solo. v.1: Conceitos básicos. Viçosa: SBCS, 2009. 595p. 5. MELO, V.F.; ALLEONI, L.R.F. (Eds.). Química e mineralogia do solo. v.2: Aplicações. Viçosa: SBCS, 2009. 685p.4. NOVAIS, R.F.; ALVAREZ V., V.H.; BARROS, N.F.; FONTES, R.L.F.; CANTARUTTI,R.B.; NEVES, J.C.L. Fertilidade do Solo. Visconde do Rio Branco: Gráfica Suprema, 2007. 1017p.5. QUAGGIO, J. A. Acidez e calagem em solos tropicais. Instituto Agronômico. 111p. (2000).6. RAIJ, B. van; ANDRADE, J.C. de; CANTARELLA, H.; QUAGGIO, J.A. Análise química para avaliação da fertilidade de solos tropicais. Raij, B. van, Andrade, J.C. de, Cantarella, H. e Quaggio, J.A. (ed.). Campinas, Instituto Agronômico, 2001. 285p.7. SANTOS, G.A; SILVA, L.S.; CANELLAS, L.P.; CAMARGO, F.A.O. (Eds). Fundamentos da matéria orgânica do solo: ecossistemas tropicais e subtropicais. Porto Alegre: Genesis. 2a Edição. 2008. 636p.",
    2)
if (-not $ok2) { throw "Find/replace for the 'Norma de recuperação:' run failed" }

Write-Host "Done."
